$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full "after" state of the From/To/Drone table (A2:D12).
# Column A holds plain numbers; columns B/C/D hold numbers-that-look-like-text
# (stored as shared strings in the original workbook).
$data = @(
    @(10, "3", "4", "1"),
    @(9,  "4", "3", "1"),
    @(8,  "1", "8", "3"),
    @(7,  "8", "1", "3"),
    @(6,  "2", "10", "7"),
    @(5,  "10", "2", "7"),
    @(4,  "1", "9", "16"),
    @(3,  "5", "1", "16"),
    @(2,  "6", "5", "16"),
    @(1,  "7", "6", "16"),
    @(0,  "9", "7", "16")
)

$startRow = 2
$endRow = $startRow + $data.Length - 1
$newRowsStart = 7

# The new rows (7-12) need the same bordered/centered style column A already
# carries on rows 2-6 (style index "s=1"). Copy that formatting down via the
# clipboard instead of Range.Style (assigning the Style object directly does
# not transfer the underlying direct formatting in this host).
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A$($newRowsStart):A$($endRow)").PasteSpecial(-4122) | Out-Null

# Pre-format the text columns as Text so the values land as shared strings
# instead of being auto-coerced to numbers. (Note: use $(...) around the
# variables here -- "B$startRow:D$endRow" would mis-parse "$startRow:D" as
# a scope-qualified variable reference in PowerShell.)
$textRange = $ws.Range("B$($startRow):D$($endRow)")
$textRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Drop the temporary Text format again so the cells fall back to the
# workbook's default (unstyled) cellXf, matching the original look.
$textRange.Style = "Normal"
